# Gate_Review_2_Tools.xlsx — "Mobile Unit" location split
#
# The tool-shed spreadsheet used one generic "Mobile Unit" location for every
# mobile-unit item. This update gives the mobile unit its two actual stops:
#   - Row 9  (MOB:M26 Jig Saw)  -> Thomas P. Ryan Center, staffed on Mondays
#   - Row 10 (MOB:W10 Studfinder) -> Edgerton Recreation Center, staffed on
#                                    Tuesdays
#
# Both the "Home Location" (col F) and "Current Location" (col G) columns
# carry the location text for these two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: rename the shared "Mobile Unit" text in place -> Monday location.
# (Both cells already hold the literal "Mobile Unit" text, so re-assigning
# the same string value to both collapses back onto a single shared entry.)
$ws.Range("F9").Value = "Mobile Unit - Thomas P. Ryan Center (Monday)"
$ws.Range("G9").Value = "Mobile Unit - Thomas P. Ryan Center (Monday)"

# Row 10: the other mobile-unit item moves to its own, distinct string ->
# Tuesday location. G10 is entered with a leading apostrophe (as the
# original author did) which forces the text quote-prefix style on that
# cell without changing the displayed value.
$ws.Range("F10").Value = "Mobile Unit - Edgerton Recreation Center (Tuesday)"
$ws.Range("G10").Value = "'Mobile Unit - Edgerton Recreation Center (Tuesday)"

# Widen columns F/G so the longer location names are fully visible (best
# effort — exact fractional "best fit" widths depend on font metrics).
$ws.Columns("F").ColumnWidth = 38.83
$ws.Columns("G").ColumnWidth = 43.6666666666667

# Leave the selection where the editor ended up while making this change.
$ws.Range("F12").Select()
